$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "대상시군구2"
$ws.Range("A1").Value = "sig_cd"
Write-Output "done"
